$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New request rows appended below the header (row 1).
# Columns: A=ID  B=fromUser  C=toUser  D=type  E=status  F=projectID  G=newTitle  H=newSupervisor

# Row 2 - deregister project request
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "YCHERN"
$ws.Range("C2").Value = "ASFLI"
$ws.Range("D2").Value = "DEREGISTERPROJECT"
$ws.Range("E2").Value = "PENDING"
$ws.Range("F2").Value = 1

# Row 3 - deregister project request
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "YCHERN"
$ws.Range("C3").Value = "ASFLI"
$ws.Range("D3").Value = "DEREGISTERPROJECT"
$ws.Range("E3").Value = "PENDING"
$ws.Range("F3").Value = 1

# Row 4 - change title request
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "YCHERN"
$ws.Range("D4").Value = "CHANGETITLE"
$ws.Range("E4").Value = "PENDING"
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = "test"

# Row 5 - change title request
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "YCHERN"
$ws.Range("C5").Value = "ASMADHUKUMAR"
$ws.Range("D5").Value = "CHANGETITLE"
$ws.Range("E5").Value = "PENDING"
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = "yays"

$ws.Range("H5").Select()
